$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header for new "id" column
$ws.Range("D2").Value = "id"

# Sequential id values for the data rows below the header (row 3 through row 17),
# including the previously-empty row 9 which now gets its own id value.
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 2
$ws.Range("D5").Value = 3
$ws.Range("D6").Value = 4
$ws.Range("D7").Value = 5
$ws.Range("D8").Value = 6
$ws.Range("D9").Value = 7
$ws.Range("D10").Value = 8
$ws.Range("D11").Value = 9
$ws.Range("D12").Value = 10
$ws.Range("D13").Value = 11
$ws.Range("D14").Value = 12
$ws.Range("D15").Value = 13
$ws.Range("D16").Value = 14
$ws.Range("D17").Value = 15

# Update the active selection to match the authored workbook state
$ws.Range("D20").Select()
